$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2049.5
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2325

$ws.Range("H43").Value = 1106.25
$ws.Range("I43").Value = 1106.25
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1106.25
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1037.25
$ws.Range("N43").ClearContents()

$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 10000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29126

$ws.Range("H70").Value = 1780.75
$ws.Range("I70").Value = 1749.5
$ws.Range("J70").Value = 1791.1666
$ws.Range("K70").Value = 5248.5
$ws.Range("L70").Value = 5373.4998
$ws.Range("M70").Value = -4978.5
$ws.Range("N70").Value = -5913.4998

$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 10000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -85632

$ws.Range("H73").Value = 1780.75
$ws.Range("I73").Value = 1749.5
$ws.Range("J73").Value = 1791.1666
$ws.Range("K73").Value = 5248.5
$ws.Range("L73").Value = 5373.4998
$ws.Range("M73").Value = -4312.5
$ws.Range("N73").Value = -7245.4998

$ws.Range("H82").Value = 729.6667
$ws.Range("I82").Value = 694.5
$ws.Range("J82").Value = 800
$ws.Range("K82").Value = 2083.5
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -1677.5
$ws.Range("N82").Value = -3212

$ws.Range("H85").Value = 729.6667
$ws.Range("I85").Value = 694.5
$ws.Range("J85").Value = 800
$ws.Range("K85").Value = 2083.5
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -679.5
$ws.Range("N85").Value = -5208

$ws.Range("H98").Value = 3411.875
$ws.Range("I98").Value = 3186
$ws.Range("K98").Value = 3186
$ws.Range("M98").Value = -1688

$ws.Range("H122").Value = 3411.875
$ws.Range("I122").Value = 3186
$ws.Range("K122").Value = 9558
$ws.Range("M122").Value = -7108

$ws.Range("H132").Value = 4348.143
$ws.Range("I132").Value = 4069
$ws.Range("K132").Value = 12207
$ws.Range("M132").Value = -9677

$ws.Range("H137").Value = 999.5
$ws.Range("I137").Value = 749
$ws.Range("K137").Value = 2247
$ws.Range("M137").Value = 303

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 19950
$ws.Range("J63").Value = 19950
$ws.Range("L63").Value = 19950
$ws.Range("N63").Value = -21322

$ws.Range("H66").Value = 19950
$ws.Range("J66").Value = 19950
$ws.Range("L66").Value = 99750
$ws.Range("N66").Value = -106614

$ws.Range("H132").Value = 3701
$ws.Range("I132").Value = 3701
$ws.Range("K132").Value = 11103
$ws.Range("M132").Value = -8573

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1667.25
$ws.Range("I20").Value = 1667.25
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1667.25
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1420.25
$ws.Range("N20").ClearContents()

$ws.Range("H80").Value = 1204.7142
$ws.Range("I80").Value = 1283.25
$ws.Range("J80").Value = 1100
$ws.Range("K80").Value = 1283.25
$ws.Range("L80").Value = 1100
$ws.Range("M80").Value = -285.25
$ws.Range("N80").Value = -3096

$ws.Range("H83").Value = 1204.7142
$ws.Range("I83").Value = 1283.25
$ws.Range("J83").Value = 1100
$ws.Range("K83").Value = 6416.25
$ws.Range("L83").Value = 5500
$ws.Range("M83").Value = -1424.25
$ws.Range("N83").Value = -15484

$ws.Range("H96").Value = 14429.5
$ws.Range("I96").Value = 14429.5
$ws.Range("K96").Value = 14429.5
$ws.Range("M96").Value = -11683.5

$ws.Range("H99").Value = 3500
$ws.Range("I99").Value = 3500
$ws.Range("K99").Value = 3500
$ws.Range("M99").Value = -2002

$ws.Range("H106").Value = 18275
$ws.Range("J106").Value = 18275
$ws.Range("L106").Value = 18275
$ws.Range("N106").Value = -20799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3141.4
$ws.Range("I58").Value = 3158.2222
$ws.Range("K58").Value = 3158.2222
$ws.Range("M58").Value = -2955.2222

$ws.Range("H62").Value = 7584.8335
$ws.Range("I62").Value = 6503
$ws.Range("J62").Value = 8666.666999999999
$ws.Range("K62").Value = 6503
$ws.Range("L62").Value = 8666.666999999999
$ws.Range("M62").Value = -5879
$ws.Range("N62").Value = -9914.666999999999

$ws.Range("H65").Value = 7584.8335
$ws.Range("I65").Value = 6503
$ws.Range("J65").Value = 8666.666999999999
$ws.Range("K65").Value = 32515
$ws.Range("L65").Value = 43333.335
$ws.Range("M65").Value = -29395
$ws.Range("N65").Value = -49573.335

$ws.Range("H132").Value = 3064
$ws.Range("I132").Value = 2387.6667
$ws.Range("K132").Value = 7163.000100000001
$ws.Range("M132").Value = -4633.000100000001

$ws.Range("H136").Value = 3141.4
$ws.Range("I136").Value = 3158.2222
$ws.Range("K136").Value = 9474.6666
$ws.Range("M136").Value = -6924.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -6058

$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2

$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 2505000
$ws.Range("I25").Value = 2505000
$ws.Range("K25").Value = 2505000
$ws.Range("M25").Value = -2504770

$ws.Range("H68").Value = 900
$ws.Range("I68").Value = 900
$ws.Range("K68").Value = 900
$ws.Range("M68").Value = -151

$ws.Range("H71").Value = 900
$ws.Range("I71").Value = 900
$ws.Range("K71").Value = 4500
$ws.Range("M71").Value = -756

$ws.Range("H94").Value = 74066
$ws.Range("J94").Value = 74066
$ws.Range("L94").Value = 74066
$ws.Range("N94").Value = -75418

$ws.Range("H122").Value = 5510.0835
$ws.Range("I122").Value = 5510.0835
$ws.Range("K122").Value = 16530.2505
$ws.Range("M122").Value = -14080.2505

$ws.Range("H132").Value = 9999
$ws.Range("I132").Value = 9999
$ws.Range("K132").Value = 29997
$ws.Range("M132").Value = -27467

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20001000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 20001000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 20001000
$ws.Range("N4").Value = -20001226
$ws.Range("M4").ClearContents()

$ws.Range("H100").Value = 302.5
$ws.Range("I100").Value = 302
$ws.Range("K100").Value = 604
$ws.Range("M100").Value = -63

$ws.Range("H122").Value = 986.25
$ws.Range("I122").Value = 986.25
$ws.Range("K122").Value = 2958.75
$ws.Range("M122").Value = -508.75
